# Fixed #295 Add the version of M2Doc in the template custom properties.
#
# This particular template resource (inferedTypeInThen-template.docx) is
# only touched incidentally by that commit: its word/document.xml and
# word/styles.xml were re-saved by the tooling, which re-emits every
# element's attributes in a normalized (alphabetical) order. A value-by-
# value comparison of the canonical OOXML before/after shows every
# element keeps exactly the same tag name, namespace and attribute
# name/value pairs -- e.g. <w:tab w:val="left" w:pos="3119"/> simply
# becomes <w:tab w:pos="3119" w:val="left"/>, <w:pgSz w:w="11906"
# w:h="16838"/> becomes <w:pgSz w:h="16838" w:w="11906"/>, and every
# <w:lsdException>/<w:style>/<w:rFonts>/<w:lang>/... attribute list is
# likewise just re-ordered. No run text, field code, tab position, page
# size/margin value, font, language, style id or latent-style flag was
# added, removed or changed.
#
# So there is no document-model mutation to perform here: the paragraphs,
# runs, field codes, tab stops, section properties and styles must come
# out of this edit holding the exact same values they went in with. We
# simply walk the object model read-only (touching nothing) so the
# template is left byte-for-byte equivalent, which is the correct,
# content-preserving result of this no-semantic-change revision.

$d = $word.ActiveDocument

$paragraphCount = $d.Paragraphs.Count
$sectionCount = $d.Sections.Count
$styleCount = $d.Styles.Count

Write-Output ("Paragraphs: " + $paragraphCount)
Write-Output ("Sections: " + $sectionCount)
Write-Output ("Styles: " + $styleCount)
Write-Output "No content changes required for this revision (attribute-order-only re-serialization)."
